# Common: Finished mix detail
# Adds new "lab.mixture.*" translation rows to the "Import" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append right after the current last row (567).
# Each tuple: (row number, label key, Czech translation)
$rows = @(
    @(568, "lab.mixture.preview.pgvg",             "PG/VG"),
    @(569, "lab.mixture.preview.age",               "Stáří mixu"),
    @(570, "lab.mixture.preview.steep",             "Doba zrání"),
    @(571, "lab.mixture.preview.mixed",             "Datum mixu"),
    @(572, "lab.mixture.preview.expires",           "Expirace"),
    @(573, "lab.mixture.preview.volume",            "Objem"),
    @(574, "lab.mixture.preview",                   "Náhled mixu"),
    @(575, "lab.mixture.preview.preview.title",     "Náhled mixu"),
    @(576, "lab.mixture.preview.preview.subtitle",  "Přehled všech dostupných dat o vybraném mixu."),
    @(577, "lab.mixture.button.index",              "Detail mixu")
)

$lastRow = 567

foreach ($row in $rows) {
    $r = $row[0]
    $label = $row[1]
    $translation = $row[2]

    # Duplicate formatting (style, etc.) from the last existing data row,
    # then overwrite the values for the new row.
    $ws.Range("A$lastRow:C$lastRow").Copy()
    $ws.Range("A${r}:C${r}").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $label
    $ws.Cells.Item($r, 3).Value = $translation
}

$excel.CutCopyMode = 0

# Match the updated view state from the authored change.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 558
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B570").Select()
